$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing data ---
# A2: "Trisna Nugraha" -> "John Doe"
$ws.Range("A2").Value = "John Doe"

# --- New header cells ---
$ws.Range("D1").Value = "Alamat"
$ws.Range("E1").Value = "No. HP"
$ws.Range("F1").Value = "Angkatan"

# --- New data cells ---
$ws.Range("D2").Value = "Jalan Maju Merdeka"
$ws.Range("E2").Value = "08123456789"
$ws.Range("F2").Value = "79/WTP"

# --- Formatting: header row (row 1) bold, filled, bordered, centered ---
$headerRange = $ws.Range("A1:F1")
$headerRange.Font.Bold = $true
$headerRange.Font.Size = 12
$headerRange.Interior.ThemeColor = 8
$headerRange.Interior.TintAndShade = 0.59999389629810485
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4108
$ws.Range("A1:F1").RowHeight = 30

# No. HP header cell (E1) uses text number format
$ws.Range("E1").NumberFormat = "@"

# --- Formatting: row 2 existing columns A,B,C centered with border ---
$ws.Range("A2:C2").HorizontalAlignment = -4108
$ws.Range("A2:C2").Borders.LineStyle = 1
$ws.Range("A2:C2").Borders.Weight = 2

$ws.Range("F2").HorizontalAlignment = -4108
$ws.Range("F2").Borders.LineStyle = 1
$ws.Range("F2").Borders.Weight = 2

# E2 (No. HP value) text format + centered + border
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").HorizontalAlignment = -4108
$ws.Range("E2").Borders.LineStyle = 1
$ws.Range("E2").Borders.Weight = 2

# --- Column widths ---
$ws.Range("A:A").ColumnWidth = 20.6640625
$ws.Range("B:B").ColumnWidth = 14.21875
$ws.Range("C:C").ColumnWidth = 17.44140625
$ws.Range("D:D").ColumnWidth = 22.5546875
$ws.Range("E:E").ColumnWidth = 19
$ws.Range("F:F").ColumnWidth = 16

# --- Page setup ---
$ws.PageSetup.Orientation = 1
